$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 is a new appended record. Column A is consistently blank/text
# in this sheet (see A2:A9), and column C holds numeric-looking quantities
# that are stored as text (see C2:C9 = "2222", "222", "333"). A leading
# apostrophe tells Excel to store the entry as literal text instead of
# auto-converting it to a number, matching the existing column formatting.
$ws.Range("A10").Value = "'"
$ws.Range("B10").Value = "احمد"
$ws.Range("C10").Value = "'222"
$ws.Range("D10").Value = "الصمود"
$ws.Range("E10").Value = "الرحلة 2"
$ws.Range("F10").Value = "C2"
$ws.Range("G10").Value = "NRC"
$ws.Range("H10").Value = "٠٢‏/٠٥‏/٢٠٢٥ ٠٢:١٠:٣٦ م"

